# Venezuela Primera Division base update (11-04-2024 00:31)
# Three pairs of rows had their match data (everything except the leading
# row-index column A) swapped between each other.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($Row1, $Row2) {
    # Columns B..AC hold all of the per-match data; column A is just the
    # running row index and must stay put.
    $addr1 = "B" + $Row1 + ":AC" + $Row1
    $addr2 = "B" + $Row2 + ":AC" + $Row2
    $range1 = $ws.Range($addr1)
    $range2 = $ws.Range($addr2)

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value2 = $vals2
    $range2.Value2 = $vals1
}

Swap-RowData 100 101
Swap-RowData 102 103
Swap-RowData 162 163
